# Complete Level17MoveList solutions. Update AchievementSuperEfficient to requirement of 35.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Achievements")

# AchievementSuperEfficient row (row 21): "Items Needed" (C21) becomes a formula
# 474+126 = 600, and "Estimated Level Acquired" (D21) becomes 35.
$ws.Range("C21").Formula = "=474+126"
$ws.Range("D21").Value = 35

# Update the sheet's active-cell selection to C14.
$ws.Range("C14").Select()
